$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.238.66"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "'1.865.32"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'243.50"
$ws.Range("E5").Value = "  +3.67%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "'0.2868"
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").Value = "'42.50"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "'0.06470"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'21.04"
$ws.Range("E11").Value = "  -3.06%  "
$ws.Range("D12").Value = "'0.07726"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("D13").Value = "'1.889.27"
$ws.Range("E13").Value = "  +1.14%  "
$ws.Range("D14").Value = "'94.95"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "'0.7035"
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("D16").Value = "'5.098"
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "'273.48"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("D18").Value = "'30.245.00"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'13.33"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("D20").Value = "'0.000007551"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").Value = "'1.0000"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'2.115.96"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("D24").Value = "'5.204"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").Value = "'6.124"
$ws.Range("E25").Value = "  -0.99%  "
$ws.Range("D26").Value = "'9.304"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("D27").Value = "'165.25"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").Value = "'18.93"
$ws.Range("E28").Value = "  +0.36%  "
$ws.Range("D29").Value = "'1.905"
$ws.Range("E29").Value = "  -2.17%  "
$ws.Range("D30").Value = "'1.371"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'0.09848"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").Value = "'1.510"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").Value = "'4.242"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").Value = "'4.019"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").Value = "'0.04728"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "'1.118"
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "'0.6908"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Value = "'2.707"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("D39").Value = "'0.01840"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").Value = "'2.741"
$ws.Range("E40").Value = "  -2.23%  "
$ws.Range("D41").Value = "'6.329"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'70.04"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "'0.8409"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'1.893"
$ws.Range("E45").Value = "  -3.23%  "
$ws.Range("D46").Value = "'0.4079"
$ws.Range("E46").Value = "  -2.36%  "
$ws.Range("D47").Value = "'101.81"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("D48").Value = "'9.236"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'7.057"
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("D50").Value = "'920.63"
$ws.Range("E50").Value = "  -4.80%  "
$ws.Range("D51").Value = "'34.88"
$ws.Range("E51").Value = "  +0.99%  "
